$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (prices) that must stay
# text cells (t="inlineStr"/"s" in the OOXML, no numeric coercion), matching
# the source data. Force text format, assign, then restore default style so
# no stray number format lingers on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.022.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.12%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4272"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +11.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3518"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07439"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.265"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.293"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.814.47"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001086"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06681"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.419"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.041.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.390"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.474"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.020.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.302"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.56"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.967"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09212"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02368"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6719"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.57%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.241"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06274"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2172"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.495"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.214"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.113"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.15%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.870"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6127"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.24%  "

$ws.Range("E48").Value = "  -3.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.049"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("E50").Value = "  -2.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07109"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.32%  "
